$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the top of the data (row 3),
# pushing every existing record (previously rows 3-86) down by one row
# (to rows 4-87). Insert a new row at position 3 to achieve that shift,
# then populate it with the new record's data.
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Vega Modelo de Temuco"
$ws.Range("C3").Value = "La Araucanía"
$ws.Range("D3").Value = 44922
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = 100112030
$ws.Range("G3").Value = "Poroto granado"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 43000
$ws.Range("L3").Value = 43000
$ws.Range("M3").Value = 43000
$ws.Range("N3").Value = "`$/saco 25 kilos"
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 1720
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
